# 2024-10-13 @ 08:36 - v4.K.3.xlsb
# Reproduce, on this generated invoice workbook, the effects of upgrading
# the GC Fiscalite invoicing macro (v4.K.3.xlsb):
#   * modFAC_Finale.bas - "Meilleur controle du format de papier (lettre)"
#       -> the print scaling baked into each invoice tab is recomputed
#   * the "2024-09-03 - 24-24481" tab is regenerated one day later, so the
#     date stamp drops out of the tab name and the invoice date label
#     advances from the 3rd to the 4th of September
#   * the "Frais d'expert en taxes" Heures/Taux cells lose their house
#     brown font colour in favour of plain black

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. The "24-24481" invoice tab: rename + bump the invoice date + font
#    colour fix on the expert-fee hours/rate line.
# ---------------------------------------------------------------------
$wsInvoice = $null
try { $wsInvoice = $wb.Worksheets.Item("2024-09-03 - 24-24481") } catch { $wsInvoice = $null }
if ($wsInvoice -eq $null) {
    try { $wsInvoice = $wb.Worksheets.Item(" - 24-24481") } catch { $wsInvoice = $null }
}

if ($wsInvoice -ne $null) {
    $wsInvoice.Name = " - 24-24481"
    $wsInvoice.Range("B21").Value = "Le 4 SEPTEMBRE 2024"
    $wsInvoice.Range("C66:D66").Font.Color = 0
}

# ---------------------------------------------------------------------
# 2. Print-scale adjustments on the four dated invoice tabs.
# ---------------------------------------------------------------------
foreach ($sheetName in @("04-10-23", "28-03-24", "11-05-24", "20-08-24")) {
    $ws = $null
    try { $ws = $wb.Worksheets.Item($sheetName) } catch { $ws = $null }
    if ($ws -ne $null) {
        $ws.PageSetup.Zoom = 62
    }
}

# ---------------------------------------------------------------------
# 3. Print-scale adjustment on the "Activites" reference tab.
# ---------------------------------------------------------------------
$wsActivites = $null
try { $wsActivites = $wb.Worksheets.Item("Activités") } catch { $wsActivites = $null }
if ($wsActivites -ne $null) {
    $wsActivites.PageSetup.Zoom = 74
}
